$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.869.03"
$ws.Range("E2").Value = "  +6.02%  "

$ws.Range("D3").Value = "2.755.41"
$ws.Range("E3").Value = "  +4.65%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "117.51"
$ws.Range("E5").Value = "  +6.22%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "331.70"
$ws.Range("E6").Value = "  +2.83%  "

$ws.Range("E7").Value = "  +2.60%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.577"
$ws.Range("E9").Value = "  +6.79%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.66"
$ws.Range("E10").Value = "  +5.32%  "

$ws.Range("E11").Value = "  +2.89%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.02"
$ws.Range("E12").Value = "  +1.16%  "

$ws.Range("E13").Value = "  +2.81%  "

$ws.Range("E14").Value = "  +5.72%  "

$ws.Range("D15").Value = "3.182.43"
$ws.Range("E15").Value = "  +4.78%  "

$ws.Range("D16").Value = "2.778.43"
$ws.Range("E16").Value = "  +5.46%  "

$ws.Range("E17").Value = "  +2.56%  "

$ws.Range("D18").Value = "51.703.45"
$ws.Range("E18").Value = "  +5.87%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.52"
$ws.Range("E19").Value = "  +5.33%  "

$ws.Range("E20").Value = "  +4.75%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.83"
$ws.Range("E21").Value = "  +2.53%  "

$ws.Range("D22").Value = "0.0₃0963"
$ws.Range("E22").Value = "  +2.24%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "278.96"
$ws.Range("E23").Value = "  +2.87%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.74"
$ws.Range("E24").Value = "  +0.78%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.65"
$ws.Range("E25").Value = "  +4.51%  "

$ws.Range("E26").Value = "  +2.81%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.14"
$ws.Range("E27").Value = "  +0.32%  "

$ws.Range("E28").Value = "  -0.02%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.24"
$ws.Range("E29").Value = "  +1.43%  "

$ws.Range("E30").Value = "  -0.03%  "

$ws.Range("E31").Value = "  +1.82%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "35.01"
$ws.Range("E32").Value = "  -0.06%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "50.45"
$ws.Range("E33").Value = "  +2.10%  "

$ws.Range("E34").Value = "  +3.02%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0822"
$ws.Range("E35").Value = "  +3.40%  "

$ws.Range("E36").Value = "  -0.10%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.06"
$ws.Range("E37").Value = "  -0.76%  "

$ws.Range("E38").Value = "  +3.08%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.97"
$ws.Range("E39").Value = "  +1.18%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.19"
$ws.Range("E40").Value = "  +1.74%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "131.00"
$ws.Range("E41").Value = "  +5.11%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0348"
$ws.Range("E42").Value = "  +11.21%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "23.17"
$ws.Range("E43").Value = "  +2.23%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.114"
$ws.Range("E44").Value = "  +2.71%  "

$ws.Range("E45").Value = "  +6.34%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.42"
$ws.Range("E46").Value = "  +14.52%  "

$ws.Range("D47").Value = "2.107.41"
$ws.Range("E47").Value = "  +2.02%  "

$ws.Range("E48").Value = "  +3.70%  "

$ws.Range("E49").Value = "  +2.30%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.57"
$ws.Range("E50").Value = "  +7.90%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.96"
$ws.Range("E51").Value = "  +0.32%  "
